$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.358.63'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '2.660.45'
$ws.Range('E3').Value = '  +1.30%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.79'
$ws.Range('E5').Value = '  -1.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.81'
$ws.Range('E6').Value = '  -2.35%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.524'
$ws.Range('E8').Value = '  -0.34%  '
$ws.Range('D9').Value = '2.659.52'
$ws.Range('E9').Value = '  +1.31%  '
$ws.Range('E10').Value = '  -2.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.169'
$ws.Range('E11').Value = '  +2.11%  '
$ws.Range('E12').Value = '  +1.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.00'
$ws.Range('E13').Value = '  -0.81%  '
$ws.Range('D14').Value = '3.151.13'
$ws.Range('E14').Value = '  +1.49%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000185'
$ws.Range('E15').Value = '  -1.82%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '72.282.30'
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.36'
$ws.Range('E17').Value = '  -1.12%  '
$ws.Range('D18').Value = '2.662.97'
$ws.Range('E18').Value = '  +1.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.26'
$ws.Range('E19').Value = '  +5.83%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.29'
$ws.Range('E20').Value = '  +3.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '370.01'
$ws.Range('E21').Value = '  -2.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.18'
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('E23').Value = '  +1.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '72.00'
$ws.Range('E24').Value = '  -1.63%  '
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.33'
$ws.Range('E26').Value = '  -1.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.77'
$ws.Range('E27').Value = '  -2.08%  '
$ws.Range('D28').Value = '2.802.57'
$ws.Range('E28').Value = '  +1.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').Value = '0.0₃0975'
$ws.Range('E30').Value = '  +1.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.13'
$ws.Range('E31').Value = '  +0.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '501.18'
$ws.Range('E32').Value = '  -3.69%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.30'
$ws.Range('E33').Value = '  -2.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.82'
$ws.Range('E34').Value = '  -0.48%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '162.85'
$ws.Range('E36').Value = '  -1.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.53'
$ws.Range('E37').Value = '  +0.95%  '
$ws.Range('B38').Value = 'WhiteBITCoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.93'
$ws.Range('E38').Value = '  -0.79%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.111'
$ws.Range('E39').Value = '  +1.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.38'
$ws.Range('E40').Value = '  -1.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.78'
$ws.Range('E41').Value = '  -2.96%  '
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.00'
$ws.Range('E43').Value = '  -1.50%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.59'
$ws.Range('E44').Value = '  -0.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.333'
$ws.Range('E45').Value = '  +0.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '156.34'
$ws.Range('E46').Value = '  +4.20%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '39.48'
$ws.Range('E47').Value = '  -0.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.75'
$ws.Range('E48').Value = '  +1.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.559'
$ws.Range('E49').Value = '  +2.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.73'
$ws.Range('E50').Value = '  +1.87%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0258'
$ws.Range('E51').Value = '  -2.03%  '
